$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the rand_digit column (J) values that changed due to re-randomized allocation
$ws.Range("J2").Value = 3
$ws.Range("J3").Value = 1
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 8
$ws.Range("J9").Value = 7
$ws.Range("J10").Value = 1
$ws.Range("J12").Value = 3
$ws.Range("J13").Value = 5
$ws.Range("J14").Value = 7
$ws.Range("J15").Value = 2
$ws.Range("J16").Value = 7
$ws.Range("J17").Value = 4
$ws.Range("J19").Value = 7
$ws.Range("J20").Value = 3
$ws.Range("J21").Value = 3
$ws.Range("J22").Value = 2
$ws.Range("J23").Value = 7
$ws.Range("J24").Value = 5
$ws.Range("J25").Value = 1
$ws.Range("J27").Value = 2
$ws.Range("J28").Value = 1
$ws.Range("J29").Value = 6
$ws.Range("J30").Value = 4
$ws.Range("J31").Value = 8
$ws.Range("J32").Value = 7
$ws.Range("J33").Value = 6
$ws.Range("J34").Value = 2
$ws.Range("J35").Value = 4
$ws.Range("J36").Value = 3
$ws.Range("J37").Value = 7
$ws.Range("J38").Value = 1
$ws.Range("J39").Value = 2
$ws.Range("J41").Value = 5
$ws.Range("J43").Value = 7
$ws.Range("J45").Value = 2
$ws.Range("J46").Value = 4
$ws.Range("J47").Value = 4
$ws.Range("J48").Value = 3
$ws.Range("J49").Value = 4
$ws.Range("J52").Value = 5
$ws.Range("J53").Value = 3
$ws.Range("J55").Value = 2
$ws.Range("J58").Value = 3
$ws.Range("J59").Value = 5
$ws.Range("J60").Value = 8
$ws.Range("J61").Value = 6
$ws.Range("J62").Value = 4
$ws.Range("J63").Value = 1
$ws.Range("J64").Value = 6
$ws.Range("J65").Value = 4
$ws.Range("J66").Value = 7
$ws.Range("J67").Value = 5
$ws.Range("J68").Value = 8
$ws.Range("J69").Value = 2
$ws.Range("J70").Value = 4
$ws.Range("J71").Value = 8
$ws.Range("J72").Value = 7
$ws.Range("J73").Value = 6
$ws.Range("J76").Value = 6
$ws.Range("J77").Value = 8
$ws.Range("J78").Value = 8
$ws.Range("J79").Value = 3
$ws.Range("J80").Value = 4
$ws.Range("J81").Value = 3
